$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) QUALITY ASSURANCE team: drop the "emelia / auditor / Nitego / Senior QA
#    Tester" entries that used to follow "Release Coordinator" (Jazz).
#    Locate by text so the script is resilient to any earlier paragraph
#    renumbering, then delete the whole run of paragraphs in one go.
# ---------------------------------------------------------------------------
$startPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "emelia") {
        $startPar = $i
        break
    }
}
if ($startPar -ne $null) {
    $endPar = $startPar + 3
    $rng = $d.Range($d.Paragraphs.Item($startPar).Range.Start, $d.Paragraphs.Item($endPar).Range.End)
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# 2) HIVE TEAM: OUTREACH: drop the "Outreach Support / Carlos Santiago /
#    Outreach Support / Emilio" entries that used to follow "Semptly".
# ---------------------------------------------------------------------------
$startPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Carlos Santiago") {
        $startPar = $i - 1
        break
    }
}
if ($startPar -ne $null) {
    $endPar = $startPar + 3
    $rng = $d.Range($d.Paragraphs.Item($startPar).Range.Start, $d.Paragraphs.Item($endPar).Range.End)
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# 3) HIVE TEAM: OUTREACH: replace the "Mark Hakkarinen / Outreach Ambassador /
#    Editor of "Your Week in SmartCash". / Email" block with the new
#    "LilyDaVine / Outreach Support" translated entry.
# ---------------------------------------------------------------------------
$startPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Mark Hakkarinen") {
        $startPar = $i
        break
    }
}
if ($startPar -ne $null) {
    $endPar = $startPar + 3
    $rng = $d.Range($d.Paragraphs.Item($startPar).Range.Start, $d.Paragraphs.Item($endPar).Range.End)
    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="105" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Exo" w:hAnsi="Exo"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="343434"/>
      <w:sz w:val="37"/>
      <w:szCs w:val="37"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Exo" w:hAnsi="Exo"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="343434"/>
      <w:sz w:val="37"/>
      <w:szCs w:val="37"/>
    </w:rPr>
    <w:t>LilyDaVine</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>
      <w:color w:val="3B3B3B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t xml:space="preserve">Outreach Support</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rng.InsertXML($xml)
}

Write-Output "done"
